$d = $word.ActiveDocument

function Split-AfterUtilities([string]$searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    $pos = $rng.End

    # $pos and $pos+1 are the two consecutive manual line breaks (<w:br/>) that
    # follow the utilities value. Turn the first one into the paragraph boundary
    # (the second one remains as the manual break starting the new paragraph).
    $d.Range($pos + 1, $pos + 1).InsertParagraphAfter()
    $d.Range($pos, $pos + 1).Text = ""
}

Split-AfterUtilities("208V/1PH;3.8A; 3/8” CW; 3/4” IW")
Split-AfterUtilities("1/2” CW, (2) 1” IW")
